# Add replanting rows (42-43) to the estimates sheet, and tidy up the
# leftover duplicate cell-formats that Excel collapses away on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New data: fruit_labor_replanting_mean_h / fruit_labor_replanting_var
#    Rows 42 and 43 already carry placeholder borders on A/B/E/F from the
#    template, so only the values need to be written.
# ---------------------------------------------------------------------

$ws.Range("A42").Value = "fruit_labor_replanting_mean_h"
$ws.Range("B42").Value = "posnorm"
$ws.Range("C42").Value = 0.3
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = 0.75
$ws.Range("F42").Value = "h"

$ws.Range("A43").Value = "fruit_labor_replanting_var"
$ws.Range("B43").Value = "tnorm_0_1"
$ws.Range("C43").Value = 0.2
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = 0.7
# F43 stays empty (it only keeps its existing border style), matching the template.

# ---------------------------------------------------------------------
# 2) Re-apply the border formatting on the existing rows so Excel's
#    save-time style cleanup collapses the now-duplicate "applyFill"
#    variants back onto the plain border styles.
#
#    A single "seed" cell is built per target style and then stamped
#    (Copy / PasteSpecial formats) onto the remaining cells that need the
#    exact same formatting - this re-uses the one resulting style index
#    instead of constructing (and registering) a fresh one per cell.
# ---------------------------------------------------------------------

$xlPasteFormats = -4122

# --- style: left + right thin border ("distribution" / "unit" columns) ---
$seedLR = $ws.Range("B10")
$seedLR.ClearFormats()
$seedLR.Borders.Item(7).Weight = 2
$seedLR.Borders.Item(7).ColorIndex = 1
$seedLR.Borders.Item(7).LineStyle = 1
$seedLR.Borders.Item(10).Weight = 2
$seedLR.Borders.Item(10).ColorIndex = 1
$seedLR.Borders.Item(10).LineStyle = 1

$leftRightBorderCells = @(
  "F10","B11","F11","B12","F12","B13","F13","B17","F17",
  "B24","F24","B25","F25","B26","F26","B27","F27","B28","F28","B29","F29",
  "B36","F36","B37","F37","B38","F38","B39","F39","B40","F40"
)
$seedLR.Copy()
foreach ($addr in $leftRightBorderCells) {
  $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# --- style: right-only thin border ("variable" / "upper" columns) ---
$seedR = $ws.Range("E10")
$seedR.ClearFormats()
$seedR.Borders.Item(10).Weight = 2
$seedR.Borders.Item(10).ColorIndex = 1
$seedR.Borders.Item(10).LineStyle = 1

$rightBorderCells = @("E17","E24","A34","A35","A36","E36","E37","E38","E39","E40")
$seedR.Copy()
foreach ($addr in $rightBorderCells) {
  $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# --- style: no special border/fill at all ---
$plainCells = @("G11","G12","G13","G20","G26","G27","G36")
foreach ($addr in $plainCells) {
  $ws.Range($addr).ClearFormats()
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) View state: the sheet had scrolled to row 10 with F42 selected; the
#    saved workbook instead keeps the selection on A10 (and Excel drops
#    the now-stale topLeftCell on save).
# ---------------------------------------------------------------------

$ws.Range("A10").Select()
